$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New roster data (player, position, team) replacing the old A2:C19 block.
$data = @(
    @("D'Angelo Russell", "PG", "Los Angeles Lakers"),
    @("Cade Cunningham", "PG,SG", "Detroit Pistons"),
    @("Derrick White", "PG,SG", "Boston Celtics"),
    @("Damian Lillard", "PG", "Milwaukee Bucks"),
    @("Cameron Johnson", "SF,PF", "Brooklyn Nets"),
    @("Jaden McDaniels", "SF,PF", "Minnesota Timberwolves"),
    @("Anthony Davis", "PF,C", "Los Angeles Lakers"),
    @("Julius Randle", "PF,C", "Minnesota Timberwolves"),
    @("Bam Adebayo", "C", "Miami Heat"),
    @("Isaiah Hartenstein", "C", "Oklahoma City Thunder"),
    @("Ben Simmons", "PG,C", "Brooklyn Nets"),
    @("Duncan Robinson", "SG,SF", "Miami Heat"),
    @("Herbert Jones", "SF,PF", "New Orleans Pelicans"),
    @("Malik Monk", "PG,SG,SF", "Sacramento Kings"),
    @("Bilal Coulibaly", "SG,SF", "Washington Wizards"),
    @("Brandon Miller", "SG,SF,PF", "Charlotte Hornets"),
    @("Brandon Ingram", "SG,SF,PF", "New Orleans Pelicans"),
    @("LaMelo Ball", "PG,SG", "Charlotte Hornets")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
